$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "column1"
$ws.Range("F1").Value = "COLUMN1"

$ws.Range("E2").Value = "a"
$ws.Range("E3").Value = "b"
$ws.Range("E4").Value = "c"
$ws.Range("E5").Value = "d"

$ws.Range("F2").Value = "e"
$ws.Range("F3").Value = "f"
$ws.Range("F4").Value = "g"
$ws.Range("F5").Value = "h"

[void]$ws.Range("F5").Select()
